# Tutorial 6 solution update for 2001ME25.xlsx
# 1) Reformat the attendance dates in column A from DD/MM/YYYY to DD-MM-YYYY
#    (kept as plain text, not converted into Excel date serials).
# 2) Refresh the Total/Real/Duplicate/Invalid/Absent tally columns (D-H)
#    for the rows whose attendance counts changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = [ordered]@{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($r in $dates.Keys) {
    $cell = $ws.Cells.Item($r, 1)
    # Force text interpretation so values like "01-08-2022" are not
    # auto-converted into date serial numbers, then restore the default
    # (unstyled) cell format so no stray number format sticks around.
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$r]
    $cell.Style = "Normal"
}

# Total Attendance Count (D), Real (E), Duplicate (F), Invalid (G), Absent (H)
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 1

$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0

$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0

$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0

$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0

$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0
